# Updated symbol list on Sun Dec 25 04:55:51 UTC 2022 with GitHub Actions
#
# The "Price" column (D) stores its values as TEXT (not numbers), and a
# couple of the "Volume(1h)" column (E) labels changed too. Because the
# COM layer auto-converts numeric-looking strings into real numbers when
# assigned directly, each numeric-looking price is written with a leading
# apostrophe to force text entry (exactly like typing '245.11 into Excel),
# and the resulting "quote prefix" cell style is reset back to Normal so
# the cell ends up as a plain, unstyled text cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$Address,
        [string]$Value
    )
    $rng = $ws.Range($Address)
    $rng.Value = "'" + $Value
    $rng.Style = "Normal"
}

Set-TextValue "D2"  "245.11"
Set-TextValue "D3"  "23.04"
Set-TextValue "D4"  "5.413"
Set-TextValue "D5"  "0.06045"
Set-TextValue "D6"  "3.394"
Set-TextValue "D7"  "0.8076"
Set-TextValue "D8"  "0.9329"
Set-TextValue "D9"  "0.1423"
Set-TextValue "D10" "0.07438"
Set-TextValue "D11" "0.03355"
Set-TextValue "D12" "0.03068"
Set-TextValue "D13" "0.09364"
Set-TextValue "D14" "3.936"
Set-TextValue "D15" "0.001594"
Set-TextValue "D16" "0.04820"
Set-TextValue "D17" "0.0005944"
Set-TextValue "D18" "0.005378"
Set-TextValue "D19" "0.004165"
Set-TextValue "D20" "0.0009859"
Set-TextValue "D21" "0.00008705"
Set-TextValue "D23" "6.441"
Set-TextValue "D40" "0.03977"

Set-TextValue "D41" "0.006404"
$ws.Range("E41").Value = "40KickTokenKICKBestin24h"

Set-TextValue "D42" "0.1073"

Set-TextValue "D44" "0.005959"
$ws.Range("E44").Value = "43LocalTradersLCT"

Set-TextValue "D45" "0.00005178"
Set-TextValue "D46" "0.00000000750"
Set-TextValue "D47" "0.0005804"
Set-TextValue "D48" "0.8505"
Set-TextValue "D49" "0.002180"
